$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1249.1111
$ws.Range("I19").Value = 1147.75
$ws.Range("J19").Value = 1330.2
$ws.Range("K19").Value = 1147.75
$ws.Range("L19").Value = 1330.2
$ws.Range("M19").Value = -972.75
$ws.Range("N19").Value = -1680.2
$ws.Range("H32").Value = 939.8
$ws.Range("J32").Value = 999.6667
$ws.Range("L32").Value = 999.6667
$ws.Range("N32").Value = -1651.6667
$ws.Range("H33").Value = 434.96155
$ws.Range("I33").Value = 260.16666
$ws.Range("K33").Value = 260.16666
$ws.Range("M33").Value = -31.16665999999998
$ws.Range("H49").Value = 2434.5
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H106").Value = 5558408
$ws.Range("I106").Value = 6063563
$ws.Range("K106").Value = 6063563
$ws.Range("M106").Value = -6062932
$ws.Range("H126").Value = 139996
$ws.Range("J126").Value = 139996
$ws.Range("L126").Value = 139996
$ws.Range("N126").Value = -149876
$ws.Range("H128").Value = 141957.5
$ws.Range("J128").Value = 141957.5
$ws.Range("L128").Value = 141957.5
$ws.Range("N128").Value = -151917.5
$ws.Range("H131").Value = 1942.6316
$ws.Range("I131").Value = 1000.93335
$ws.Range("J131").Value = 5474
$ws.Range("K131").Value = 3002.80005
$ws.Range("L131").Value = 16422
$ws.Range("M131").Value = 2037.19995
$ws.Range("N131").Value = -26502
$ws.Range("H132").Value = 1680.1
$ws.Range("I132").Value = 1594.2222
$ws.Range("K132").Value = 4782.6666
$ws.Range("M132").Value = -2252.6666
$ws.Range("H135").Value = 1351.6364
$ws.Range("I135").Value = 1464.2354
$ws.Range("J135").Value = 968.8
$ws.Range("K135").Value = 13178.1186
$ws.Range("L135").Value = 8719.199999999999
$ws.Range("M135").Value = -10643.1186
$ws.Range("N135").Value = -13789.2
$ws.Range("H138").Value = 2718.192
$ws.Range("I138").Value = 1861.4445
$ws.Range("J138").Value = 2908.5803
$ws.Range("K138").Value = 5584.333500000001
$ws.Range("L138").Value = 8725.740900000001
$ws.Range("M138").Value = -444.3335000000006
$ws.Range("N138").Value = -19005.7409
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23811078
$ws.Range("I32").Value = 25642308
$ws.Range("K32").Value = 25642308
$ws.Range("M32").Value = -25642021
$ws.Range("H45").Value = 4701.5713
$ws.Range("I45").Value = 5439.6
$ws.Range("J45").Value = 4291.5557
$ws.Range("K45").Value = 5439.6
$ws.Range("L45").Value = 4291.5557
$ws.Range("M45").Value = -5062.6
$ws.Range("N45").Value = -5045.5557
$ws.Range("H61").Value = 8021.533
$ws.Range("I61").Value = 13933.333
$ws.Range("J61").Value = 4080.3333
$ws.Range("K61").Value = 13933.333
$ws.Range("L61").Value = 4080.3333
$ws.Range("M61").Value = -13721.333
$ws.Range("N61").Value = -4504.3333
$ws.Range("H120").Value = 67049.5
$ws.Range("J120").Value = 67049.5
$ws.Range("L120").Value = 67049.5
$ws.Range("N120").Value = -76725.5
$ws.Range("H132").Value = 4866.773
$ws.Range("I132").Value = 5921.364
$ws.Range("K132").Value = 17764.092
$ws.Range("M132").Value = -15234.092
$ws.Range("H136").Value = 8021.533
$ws.Range("I136").Value = 13933.333
$ws.Range("J136").Value = 4080.3333
$ws.Range("K136").Value = 41799.999
$ws.Range("L136").Value = 12240.9999
$ws.Range("M136").Value = -39249.999
$ws.Range("N136").Value = -17340.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 40550.152
$ws.Range("I20").Value = 57461.332
$ws.Range("K20").Value = 57461.332
$ws.Range("M20").Value = -57214.332
$ws.Range("H134").Value = 3570.2703
$ws.Range("I134").Value = 3041.111
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 9123.332999999999
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -6588.332999999999
$ws.Range("N134").Value = -20067
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3924.7827
$ws.Range("J31").Value = 4359
$ws.Range("L31").Value = 4359
$ws.Range("N31").Value = -4949
$ws.Range("H34").Value = 3924.7827
$ws.Range("J34").Value = 4359
$ws.Range("L34").Value = 4359
$ws.Range("N34").Value = -4763
$ws.Range("H58").Value = 3368.4666
$ws.Range("I58").Value = 2869.1667
$ws.Range("J58").Value = 4367.067
$ws.Range("K58").Value = 2869.1667
$ws.Range("L58").Value = 4367.067
$ws.Range("M58").Value = -2666.1667
$ws.Range("N58").Value = -4773.067
$ws.Range("H70").Value = 41999.668
$ws.Range("J70").Value = 41999.668
$ws.Range("L70").Value = 41999.668
$ws.Range("N70").Value = -42629.668
$ws.Range("H73").Value = 41999.668
$ws.Range("J73").Value = 41999.668
$ws.Range("L73").Value = 41999.668
$ws.Range("N73").Value = -44183.668
$ws.Range("H122").Value = 5987.8
$ws.Range("I122").Value = 4978.6
$ws.Range("J122").Value = 6997
$ws.Range("K122").Value = 14935.8
$ws.Range("L122").Value = 20991
$ws.Range("M122").Value = -12485.8
$ws.Range("N122").Value = -25891
$ws.Range("H136").Value = 3368.4666
$ws.Range("I136").Value = 2869.1667
$ws.Range("J136").Value = 4367.067
$ws.Range("K136").Value = 8607.500100000001
$ws.Range("L136").Value = 13101.201
$ws.Range("M136").Value = -6057.500100000001
$ws.Range("N136").Value = -18201.201
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2220
$ws.Range("J21").Value = 2480
$ws.Range("L21").Value = 7440
$ws.Range("N21").Value = -7786
$ws.Range("H23").Value = 62.25
$ws.Range("J23").Value = 76.333336
$ws.Range("L23").Value = 229.000008
$ws.Range("N23").Value = -699.000008
$ws.Range("H42").Value = 16666.666
$ws.Range("J42").Value = 16666.666
$ws.Range("L42").Value = 49999.99800000001
$ws.Range("N42").Value = -51067.99800000001
$ws.Range("H128").Value = 1979899
$ws.Range("I128").Value = 1979899
$ws.Range("K128").Value = 5939697
$ws.Range("M128").Value = -5934717
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 410.41666
$ws.Range("I2").Value = 375
$ws.Range("J2").Value = 587.5
$ws.Range("K2").Value = 375
$ws.Range("L2").Value = 587.5
$ws.Range("M2").Value = -262
$ws.Range("N2").Value = -813.5
$ws.Range("H97").Value = 936.8
$ws.Range("J97").Value = 2065.25
$ws.Range("L97").Value = 2065.25
$ws.Range("N97").Value = -3057.25
$ws.Range("H109").Value = 104994.5
$ws.Range("J109").Value = 104994.5
$ws.Range("L109").Value = 104994.5
$ws.Range("N109").Value = -107074.5
$ws.Range("H113").Value = 24472.5
$ws.Range("I113").Value = 11254.409
$ws.Range("K113").Value = 11254.409
$ws.Range("M113").Value = -9084.409
$ws.Range("H122").Value = 3154.3333
$ws.Range("I122").Value = 3334.2
$ws.Range("J122").Value = 2255
$ws.Range("K122").Value = 10002.6
$ws.Range("L122").Value = 6765
$ws.Range("M122").Value = -7552.599999999999
$ws.Range("N122").Value = -11665
$ws.Range("H126").Value = 2728.7778
$ws.Range("I126").Value = 2489.8572
$ws.Range("K126").Value = 7469.571599999999
$ws.Range("M126").Value = -4999.571599999999
$ws.Range("H132").Value = 4438.875
$ws.Range("I132").Value = 4602.4
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 13807.2
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -11277.2
$ws.Range("N132").Value = -17559.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3300.2273
$ws.Range("I7").Value = 3178.5881
$ws.Range("J7").Value = 3713.8
$ws.Range("K7").Value = 3178.5881
$ws.Range("L7").Value = 3713.8
$ws.Range("M7").Value = -3066.5881
$ws.Range("N7").Value = -3937.8
$ws.Range("H61").Value = 4179.9
$ws.Range("I61").Value = 1682.3334
$ws.Range("K61").Value = 1682.3334
$ws.Range("M61").Value = -1480.3334
$ws.Range("H113").Value = 4179.9
$ws.Range("I113").Value = 1682.3334
$ws.Range("K113").Value = 1682.3334
$ws.Range("M113").Value = 487.6666
$ws.Range("H126").Value = 3300.2273
$ws.Range("I126").Value = 3178.5881
$ws.Range("J126").Value = 3713.8
$ws.Range("K126").Value = 9535.764299999999
$ws.Range("L126").Value = 11141.4
$ws.Range("M126").Value = -7065.764299999999
$ws.Range("N126").Value = -16081.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47302.25
$ws.Range("J70").Value = 38736.668
$ws.Range("L70").Value = 38736.668
$ws.Range("N70").Value = -39366.668
$ws.Range("H73").Value = 47302.25
$ws.Range("J73").Value = 38736.668
$ws.Range("L73").Value = 38736.668
$ws.Range("N73").Value = -40920.668
$ws.Range("H75").Value = 114499.5
$ws.Range("J75").Value = 114499.5
$ws.Range("L75").Value = 114499.5
$ws.Range("N75").Value = -116371.5
$ws.Range("H78").Value = 114499.5
$ws.Range("J78").Value = 114499.5
$ws.Range("L78").Value = 343498.5
$ws.Range("N78").Value = -352858.5
$ws.Range("H126").Value = 2715.5715
$ws.Range("I126").Value = 2834.6667
$ws.Range("J126").Value = 2001
$ws.Range("K126").Value = 8504.000100000001
$ws.Range("L126").Value = 6003
$ws.Range("M126").Value = -6034.000100000001
$ws.Range("N126").Value = -10943
$ws.Range("H136").Value = 2424.5454
$ws.Range("I136").Value = 1940.1666
$ws.Range("K136").Value = 5820.4998
$ws.Range("M136").Value = -3270.4998
